$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 14 ("Fruta / hortaliza, semanal"),
# pushing all subsequent records (previously rows 14-51) down by one row
# (now rows 15-52).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record's data.
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 45246
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100112022
$ws.Cells.Item(14, 7).Value = "Arveja Verde"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 20000
$ws.Cells.Item(14, 12).Value = 20000
$ws.Cells.Item(14, 13).Value = 20000
$ws.Cells.Item(14, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 800
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
